$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1).Range
$p1.InsertParagraphAfter()
$newPara = $d.Paragraphs(2).Range
$newPara.Text = "GitHub Link: https://github.com/ABJ-Gore/student-performance-analyzer"

$urlRange = $d.Paragraphs(2).Range.Duplicate
$foundUrl = $urlRange.Find.Execute("https://github.com/ABJ-Gore/student-performance-analyzer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$urlRange.Bold = 0
$urlRange.Bold = 1

$hlRange = $d.Paragraphs(2).Range.Duplicate
$foundHl = $hlRange.Find.Execute("//github.com/ABJ-Gore/student-performance-analyzer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hlRange.Bold = 0
$hlRange.Bold = 1

$hlRange2 = $d.Paragraphs(2).Range.Duplicate
$foundHl2 = $hlRange2.Find.Execute("//github.com/ABJ-Gore/student-performance-analyzer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$link = $d.Hyperlinks.Add($hlRange2, "https://github.com/ABJ-Gore/student-performance-analyzer", $null, $null, "//github.com/ABJ-Gore/student-performance-analyzer")

$hlRange3 = $d.Paragraphs(2).Range.Duplicate
$foundHl3 = $hlRange3.Find.Execute("//github.com/ABJ-Gore/student-performance-analyzer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hlRange3.Bold = 1
$hlRange3.Bold = 0
$hlRange3.Bold = 1

Write-Output "done"
